$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 4, 7, 150, 4, 69400, 0.2164772727272727),
    @(0, 4, 6, 150, 4, 63200, 0.2145833333333333),
    @(0, 4, 6, 100, 4, 170200, 0.1921882832080201),
    @(0, 4, 5, 100, 4, 151800, 0.1890533625730994),
    @(0, 4, 8, 150, 4, 59400, 0.1770021645021645),
    @(0, 5, 6, 150, 4, 113800, 0.30625),
    @(0, 5, 5, 150, 4, 90200, 0.250297619047619)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}
